# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the f633edfe-... file in both the zh-cn and de-de sheets,
# plus the matching status text on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Excel's ColumnWidth property (characters) maps to the raw OOXML column
# width (roughly characters + ~0.83 in this runtime's pixel model), so use
# the inverse value that round-trips to exactly width="40" in the XML.
$errorColWidth = 39.1666666666667

# Overview sheet: row 3 is the f633edfe-... file; its Status (zh-cn / de-de
# columns, E and F) moves from "Ready for handoff" to "Handback transform failed".
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn sheet: row 3 (the f633edfe-... file) Status column also shared the
# same "Ready for handoff" string, so it now reads the same new status, and
# gets an Error Detail message in column P; that column is widened to fit it.
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = "Handback file name: xuxic4vr.wje is different with handoff file name: f633edfe-6f8b-4ed4-a778-e7abc896bf87.165c0399c9b1f112d0ffd03926ffa4faa2c2b3f9.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = $errorColWidth

# de-de sheet: same Status update, same Error Detail message (localized
# suffix) for row 3, and the same column widening.
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = "Handback file name: xuxic4vr.wje is different with handoff file name: f633edfe-6f8b-4ed4-a778-e7abc896bf87.165c0399c9b1f112d0ffd03926ffa4faa2c2b3f9.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = $errorColWidth
